$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 4.460132666666667
$ws.Range("H2").Value = 13.380398
$ws.Range("I2").Value = 0.02674725343762847
$ws.Range("J2").Value = 0.02674725343762847
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 98.946724
$ws.Range("N2").Value = 296.840172
$ws.Range("O2").Value = 0.2098009692989996
$ws.Range("P2").Value = 0.2098009692989996
$ws.Range("Q2").Value = 441.3155159720507
$ws.Range("R2").Value = 3971.839643748456
$ws.Range("S2").Value = 0.005611599697300453
$ws.Range("T2").Value = 0.005611599697300453

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 4.460132666666667
$ws.Range("H3").Value = 13.380398
$ws.Range("I3").Value = 0.02674725343762847
$ws.Range("J3").Value = 0.02674725343762847
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 163.0062356666667
$ws.Range("N3").Value = 489.018707
$ws.Range("O3").Value = 0.345629090707923
$ws.Range("P3").Value = 0.3456290907079231
$ws.Range("Q3").Value = 727.0294365672651
$ws.Range("R3").Value = 6543.264929105386
$ws.Range("S3").Value = 0.009244628884581899
$ws.Range("T3").Value = 0.009244628884581899

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 4.460132666666667
$ws.Range("H4").Value = 13.380398
$ws.Range("I4").Value = 0.02674725343762847
$ws.Range("J4").Value = 0.02674725343762847
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 65.39610666666668
$ws.Range("N4").Value = 196.18832
$ws.Range("O4").Value = 0.1386621609326595
$ws.Range("P4").Value = 0.1386621609326595
$ws.Range("Q4").Value = 291.6753116168179
$ws.Range("R4").Value = 2625.07780455136
$ws.Range("S4").Value = 0.003708831960675068
$ws.Range("T4").Value = 0.003708831960675069

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 4.460132666666667
$ws.Range("H5").Value = 13.380398
$ws.Range("I5").Value = 0.02674725343762847
$ws.Range("J5").Value = 0.02674725343762847
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 144.2727966666667
$ws.Range("N5").Value = 432.81839
$ws.Range("O5").Value = 0.3059077790604178
$ws.Range("P5").Value = 0.3059077790604179
$ws.Range("Q5").Value = 643.4758133243578
$ws.Range("R5").Value = 5791.28231991922
$ws.Range("S5").Value = 0.008182192895071052
$ws.Range("T5").Value = 0.008182192895071053

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 128.9378306666667
$ws.Range("H6").Value = 386.813492
$ws.Range("I6").Value = 0.7732354825034408
$ws.Range("J6").Value = 0.7732354825034408
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 98.946724
$ws.Range("N6").Value = 296.840172
$ws.Range("O6").Value = 0.2098009692989996
$ws.Range("P6").Value = 0.2098009692989996
$ws.Range("Q6").Value = 12757.9759441334
$ws.Range("R6").Value = 114821.7834972006
$ws.Range("S6").Value = 0.1622255537256015
$ws.Range("T6").Value = 0.1622255537256015

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 128.9378306666667
$ws.Range("H7").Value = 386.813492
$ws.Range("I7").Value = 0.7732354825034408
$ws.Range("J7").Value = 0.7732354825034408
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 163.0062356666667
$ws.Range("N7").Value = 489.018707
$ws.Range("O7").Value = 0.345629090707923
$ws.Range("P7").Value = 0.3456290907079231
$ws.Range("Q7").Value = 21017.67041199942
$ws.Range("R7").Value = 189159.0337079949
$ws.Range("S7").Value = 0.2672526767207664
$ws.Range("T7").Value = 0.2672526767207664

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 128.9378306666667
$ws.Range("H8").Value = 386.813492
$ws.Range("I8").Value = 0.7732354825034408
$ws.Range("J8").Value = 0.7732354825034408
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 65.39610666666668
$ws.Range("N8").Value = 196.18832
$ws.Range("O8").Value = 0.1386621609326595
$ws.Range("P8").Value = 0.1386621609326595
$ws.Range("Q8").Value = 8432.03212764594
$ws.Range("R8").Value = 75888.28914881345
$ws.Range("S8").Value = 0.1072185029137347
$ws.Range("T8").Value = 0.1072185029137347

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 128.9378306666667
$ws.Range("H9").Value = 386.813492
$ws.Range("I9").Value = 0.7732354825034408
$ws.Range("J9").Value = 0.7732354825034408
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 144.2727966666667
$ws.Range("N9").Value = 432.81839
$ws.Range("O9").Value = 0.3059077790604178
$ws.Range("P9").Value = 0.3059077790604179
$ws.Range("Q9").Value = 18602.2214264131
$ws.Range("R9").Value = 167419.9928377179
$ws.Range("S9").Value = 0.2365387491433381
$ws.Range("T9").Value = 0.2365387491433382

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 3.059082666666667
$ws.Range("H10").Value = 9.177248000000002
$ws.Range("I10").Value = 0.01834520752790531
$ws.Range("J10").Value = 0.01834520752790531
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 98.946724
$ws.Range("N10").Value = 296.840172
$ws.Range("O10").Value = 0.2098009692989996
$ws.Range("P10").Value = 0.2098009692989996
$ws.Range("Q10").Value = 302.6862083118507
$ws.Range("R10").Value = 2724.175874806657
$ws.Range("S10").Value = 0.003848842321345837
$ws.Range("T10").Value = 0.003848842321345837

# Row 11
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 3.059082666666667
$ws.Range("H11").Value = 9.177248000000002
$ws.Range("I11").Value = 0.01834520752790531
$ws.Range("J11").Value = 0.01834520752790531
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 163.0062356666667
$ws.Range("N11").Value = 489.018707
$ws.Range("O11").Value = 0.345629090707923
$ws.Range("P11").Value = 0.3456290907079231
$ws.Range("Q11").Value = 498.6495500864819
$ws.Range("R11").Value = 4487.845950778337
$ws.Range("S11").Value = 0.006340637396718055
$ws.Range("T11").Value = 0.006340637396718056

# Row 12
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 3.059082666666667
$ws.Range("H12").Value = 9.177248000000002
$ws.Range("I12").Value = 0.01834520752790531
$ws.Range("J12").Value = 0.01834520752790531
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 65.39610666666668
$ws.Range("N12").Value = 196.18832
$ws.Range("O12").Value = 0.1386621609326595
$ws.Range("P12").Value = 0.1386621609326595
$ws.Range("Q12").Value = 200.0520963714845
$ws.Range("R12").Value = 1800.468867343361
$ws.Range("S12").Value = 0.002543786118577441
$ws.Range("T12").Value = 0.002543786118577442

# Row 13
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 3.059082666666667
$ws.Range("H13").Value = 9.177248000000002
$ws.Range("I13").Value = 0.01834520752790531
$ws.Range("J13").Value = 0.01834520752790531
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 144.2727966666667
$ws.Range("N13").Value = 432.81839
$ws.Range("O13").Value = 0.3059077790604178
$ws.Range("P13").Value = 0.3059077790604179
$ws.Range("Q13").Value = 441.3424115545245
$ws.Range("R13").Value = 3972.081703990721
$ws.Range("S13").Value = 0.00561194169126397
$ws.Range("T13").Value = 0.005611941691263971

# Row 14
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 30.29400666666666
$ws.Range("H14").Value = 90.88201999999998
$ws.Range("I14").Value = 0.1816720565310254
$ws.Range("J14").Value = 0.1816720565310254
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 98.946724
$ws.Range("N14").Value = 296.840172
$ws.Range("O14").Value = 0.2098009692989996
$ws.Range("P14").Value = 0.2098009692989996
$ws.Range("Q14").Value = 2997.492716500826
$ws.Range("R14").Value = 26977.43444850744
$ws.Range("S14").Value = 0.03811497355475178
$ws.Range("T14").Value = 0.03811497355475178

# Row 15
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 30.29400666666666
$ws.Range("H15").Value = 90.88201999999998
$ws.Range("I15").Value = 0.1816720565310254
$ws.Range("J15").Value = 0.1816720565310254
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 163.0062356666667
$ws.Range("N15").Value = 489.018707
$ws.Range("O15").Value = 0.345629090707923
$ws.Range("P15").Value = 0.3456290907079231
$ws.Range("Q15").Value = 4938.111989994237
$ws.Range("R15").Value = 44443.00790994813
$ws.Range("S15").Value = 0.0627911477058567
$ws.Range("T15").Value = 0.06279114770585671

# Row 16
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 30.29400666666666
$ws.Range("H16").Value = 90.88201999999998
$ws.Range("I16").Value = 0.1816720565310254
$ws.Range("J16").Value = 0.1816720565310254
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 65.39610666666668
$ws.Range("N16").Value = 196.18832
$ws.Range("O16").Value = 0.1386621609326595
$ws.Range("P16").Value = 0.1386621609326595
$ws.Range("Q16").Value = 1981.110091334044
$ws.Range("R16").Value = 17829.9908220064
$ws.Range("S16").Value = 0.02519103993967225
$ws.Range("T16").Value = 0.02519103993967226

# Row 17
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 30.29400666666666
$ws.Range("H17").Value = 90.88201999999998
$ws.Range("I17").Value = 0.1816720565310254
$ws.Range("J17").Value = 0.1816720565310254
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 144.2727966666667
$ws.Range("N17").Value = 432.81839
$ws.Range("O17").Value = 0.3059077790604178
$ws.Range("P17").Value = 0.3059077790604179
$ws.Range("Q17").Value = 4370.601064038644
$ws.Range("R17").Value = 39335.4095763478
$ws.Range("S17").Value = 0.05557489533074466
$ws.Range("T17").Value = 0.05557489533074467

